# "Codes cleaning - data cleaning & PSW"
# Refresh the OR / lower / higher estimates on the "by_prov" sheet after
# re-running the cleaning pipeline. Only the numeric results move; the
# province / group labels in columns A and E stay as they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by_prov")

# Row => OR (col B), lower (col C), higher (col D)
$ws.Cells.Item(2, 2).Value  = 3.13
$ws.Cells.Item(2, 3).Value  = 1.9
$ws.Cells.Item(2, 4).Value  = 5.39

$ws.Cells.Item(3, 2).Value  = 1.4
$ws.Cells.Item(3, 3).Value  = 1.18
$ws.Cells.Item(3, 4).Value  = 1.65

$ws.Cells.Item(4, 2).Value  = 3.26
$ws.Cells.Item(4, 3).Value  = 2.16
$ws.Cells.Item(4, 4).Value  = 5.01

$ws.Cells.Item(5, 2).Value  = 1.44
$ws.Cells.Item(5, 3).Value  = 1.22
$ws.Cells.Item(5, 4).Value  = 1.7

$ws.Cells.Item(6, 2).Value  = 4.19
$ws.Cells.Item(6, 3).Value  = 2.93
$ws.Cells.Item(6, 4).Value  = 6.15

$ws.Cells.Item(7, 2).Value  = 1.44
$ws.Cells.Item(7, 3).Value  = 1.24
$ws.Cells.Item(7, 4).Value  = 1.67

$ws.Cells.Item(8, 2).Value  = 1.89
$ws.Cells.Item(8, 3).Value  = 1.34
$ws.Cells.Item(8, 4).Value  = 2.7

# Row 9 only changed the "higher" value.
$ws.Cells.Item(9, 4).Value  = 1.75

$ws.Cells.Item(10, 2).Value = 7.26
$ws.Cells.Item(10, 3).Value = 4.52
$ws.Cells.Item(10, 4).Value = 12.4

$ws.Cells.Item(11, 2).Value = 1.46
$ws.Cells.Item(11, 3).Value = 1.25
$ws.Cells.Item(11, 4).Value = 1.71
